$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 811759
$ws.Range("C2").Value = 835944.9489999991
$ws.Range("D2").Value = 959725.0510000009
$ws.Range("F2").Value = 741612.915
$ws.Range("G2").Value = 1161081.284
$ws.Range("H2").Value = 2133796
$ws.Range("I2").Value = 2976497.263
$ws.Range("J2").Value = 2134978.462
$ws.Range("K2").Value = 1362887
$ws.Range("N2").Value = 1220744
$ws.Range("O2").Value = 1246828
$ws.Range("P2").Value = 1019833
$ws.Range("R2").Value = 1204937
$ws.Range("S2").Value = 1825975.33
$ws.Range("T2").Value = 2531543.67
$ws.Range("U2").Value = 3160953
$ws.Range("V2").Value = 2820420.585999999
$ws.Range("W2").Value = 2094488.738
$ws.Range("X2").Value = 1232448
$ws.Range("Y2").Value = 997893
$ws.Range("Z2").Value = 1021834.940000001
$ws.Range("AA2").Value = 1225809.770999998
$ws.Range("AB2").Value = 1081975.289000001
$ws.Range("AD2").Value = 900861
$ws.Range("AE2").Value = 1500648
$ws.Range("AF2").Value = 2514295
$ws.Range("AG2").Value = 3493064
$ws.Range("AH2").Value = 2269219
$ws.Range("AI2").Value = 1598814.403000001
$ws.Range("AJ2").Value = 1140410
$ws.Range("AK2").Value = 810234
$ws.Range("AL2").Value = 1219717
$ws.Range("AM2").Value = 1099240
$ws.Range("AN2").Value = 1115052
$ws.Range("AP2").Value = 957949
$ws.Range("AQ2").Value = 2031447
$ws.Range("AR2").Value = 2955074
$ws.Range("AS2").Value = 3551912
$ws.Range("AT2").Value = 3226166.652000001
$ws.Range("AU2").Value = 1880328
$ws.Range("AV2").Value = 1390254
$ws.Range("AW2").Value = 964408
$ws.Range("AX2").Value = 1183197
$ws.Range("AY2").Value = 1532346
$ws.Range("AZ2").Value = 1335087
$ws.Range("BB2").Value = 1021321
